$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "neel"
$ws.Range("A2").Value = "saspara"

$ws.Range("A2").Select()
